$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 8
    3  = 9
    4  = 7
    5  = 4
    6  = 8
    7  = 14
    8  = 6
    9  = 4
    10 = 7
    11 = 7
    12 = 8
    13 = 8
    14 = 9
    15 = 8
    16 = 8
    17 = 11
    18 = 12
    19 = 7
    20 = 6
    21 = 6
    22 = 6
    23 = 8
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
